$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 and J1 (copy formatting style from an existing header cell)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New column data I2:J41
$data = @(
  @(2,3),
  @(8,9),
  @(9,9),
  @(9,9),
  @(5,5),
  @(9,9),
  @(7,7),
  @(8,8),
  @(4,5),
  @(6,6),
  @(5,6),
  @(7,8),
  @(5,6),
  @(6,6),
  @(7,7),
  @(7,7),
  @(9,9),
  @(9,9),
  @(7,9),
  @(7,8),
  @(5,5),
  @(6,6),
  @(7,8),
  @(7,8),
  @(8,9),
  @(8,8),
  @(7,7),
  @(8,9),
  @(3,4),
  @(7,9),
  @(6,7),
  @(5,6),
  @(6,6),
  @(8,8),
  @(7,7),
  @(8,8),
  @(6,6),
  @(8,8),
  @(9,9),
  @(4,4)
)

$row = 2
foreach ($pair in $data) {
  $ws.Cells.Item($row, 9).Value = $pair[0]
  $ws.Cells.Item($row, 10).Value = $pair[1]
  $row = $row + 1
}
